# Update gh-pages to output generated at 456a3b4
# Increments the "想去人数" (want-to-go count, column F) for a handful of
# events that appear on the "展览" (Exhibitions), "演出" (Performances) and
# the aggregated "全部类型" (All types) sheets.

$wb = $excel.ActiveWorkbook

$wsExhibitions = $wb.Worksheets.Item("展览")
$wsPerformances = $wb.Worksheets.Item("演出")
$wsAll = $wb.Worksheets.Item("全部类型")

# 展览 sheet updates
$wsExhibitions.Range("F7").Value = 42
$wsExhibitions.Range("F8").Value = 681
$wsExhibitions.Range("F10").Value = 173
$wsExhibitions.Range("F14").Value = 533
$wsExhibitions.Range("F15").Value = 478
$wsExhibitions.Range("F20").Value = 2565
$wsExhibitions.Range("F27").Value = 123
$wsExhibitions.Range("F29").Value = 931
$wsExhibitions.Range("F31").Value = 129

# 演出 sheet updates
$wsPerformances.Range("F9").Value = 288
$wsPerformances.Range("F13").Value = 515

# 全部类型 sheet updates (aggregate of the above, same events different rows)
$wsAll.Range("F16").Value = 42
$wsAll.Range("F17").Value = 681
$wsAll.Range("F20").Value = 173
$wsAll.Range("F23").Value = 533
$wsAll.Range("F24").Value = 478
$wsAll.Range("F28").Value = 2565
$wsAll.Range("F32").Value = 288
$wsAll.Range("F34").Value = 123
$wsAll.Range("F36").Value = 931
$wsAll.Range("F37").Value = 515
$wsAll.Range("F40").Value = 129
